$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("November")

# New headers for the additional DB-style columns (H1:P1)
$ws.Range("H1").Value = "Check Out Date"
$ws.Range("I1").Value = "Paid"
$ws.Range("J1").Value = "Payment Method"
$ws.Range("K1").Value = "Amount Paid"
$ws.Range("L1").Value = "Currency"
$ws.Range("M1").Value = "Rooms"
$ws.Range("N1").Value = "Beds"
$ws.Range("O1").Value = "Person"
$ws.Range("P1").Value = "Notes"

# Column widths for the new columns
$ws.Range("H1:I1").ColumnWidth = 10
$ws.Range("J1").ColumnWidth = 15
$ws.Range("K1:P1").ColumnWidth = 10

# Existing rows shift check-in date forward by one day
$ws.Range("A2").Value = (Get-Date -Year 2024 -Month 11 -Day 26).Date
$ws.Range("A3").Value = (Get-Date -Year 2024 -Month 11 -Day 26).Date

# Row 3 guest data corrections
$ws.Range("B3").Value = "Stef"

# D3 and G3 need to stay text (not auto-converted to numbers), so go through
# a formula -> values-only paste round trip to force a text shared-string.
$ws.Range("D3").Formula = "=TEXT(1,""0"")"
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("G3").Formula = "=TEXT(123,""0"")"
$ws.Range("G3").Copy()
$ws.Range("G3").PasteSpecial(-4163)

# New Check Out Date column (H) values, reusing the date style from column A
$ws.Range("A2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = (Get-Date -Year 2024 -Month 11 -Day 29).Date

$ws.Range("A2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = (Get-Date -Year 2024 -Month 11 -Day 29).Date
